$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 6.116799999999991
$ws.Range("A9").Value = -20.27869999999998
$ws.Range("A18").Value = -23.13100000000001
$ws.Range("A20").Value = -22.22160000000003
$ws.Range("D21").Value = -7.6832
